$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (date, D, E, F, G, H)
$rows = @{
    3  = @("28-07-2022", 1, 0, 0, 1, 1)
    4  = @("01-08-2022", 1, 1, 0, 0, 0)
    5  = @("04-08-2022", 1, 1, 0, 0, 0)
    6  = @("08-08-2022", 1, 1, 0, 0, 0)
    7  = @("11-08-2022", 0, 0, 0, 0, 1)
    8  = @("15-08-2022", 0, 0, 0, 0, 1)
    9  = @("18-08-2022", 0, 0, 0, 0, 1)
    10 = @("22-08-2022", 1, 1, 0, 0, 0)
    11 = @("25-08-2022", 1, 1, 0, 0, 0)
    12 = @("29-08-2022", 1, 1, 0, 0, 0)
    13 = @("01-09-2022", 1, 1, 0, 0, 0)
    14 = @("05-09-2022", 1, 1, 0, 0, 0)
    15 = @("08-09-2022", 0, 0, 0, 0, 1)
    16 = @("12-09-2022", 0, 0, 0, 0, 1)
    17 = @("15-09-2022", 0, 0, 0, 0, 1)
    18 = @("19-09-2022", 0, 0, 0, 0, 1)
    19 = @("22-09-2022", 0, 0, 0, 0, 1)
    20 = @("26-09-2022", 0, 0, 0, 0, 1)
    21 = @("29-09-2022", 0, 0, 0, 0, 1)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
    $ws.Cells.Item($r, 8).Value = $vals[5]
}
